# Draft before submission to Barcelona
# - Remove the "SE+AR" column (column F) entirely, including its header
#   and data cells.
# - E2 ("InfAV" row, "FIRE+SV" column) changes from the shared string
#   "N/A" to the numeric value 0.
# - InfATV row (row 4) SPF/SCE moments are recomputed: B4 0.115 -> 0.918,
#   C4 0.417 -> 3.763.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the SE+AR column (F) and shift everything left.
$ws.Columns.Item(6).Delete()

# Update the recomputed / corrected values.
$ws.Range("E2").Value = 0
$ws.Range("B4").Value = 0.918
$ws.Range("C4").Value = 3.763
